$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Christian Aliwate"
$ws.Range("B4").Value = "christian.aliwate@gmail.com"
$ws.Range("C4").Value = "super_admin@logistic"
$ws.Range("D4").Value = "l0gistic@!"
$ws.Range("E4").Value = "0000-00-00 00:00:00"
$ws.Range("F4").Value = "2023-05-03 20:52:50"
$ws.Range("G4").Value = "uploads/fleet_management/users/FMD-64525922/1683118370_christian_aliwate.jpg"
$ws.Range("H4").Value = "active"
